$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NuevoPrecio")

# Listado nuevo en lista de precios: actualizar columna "precio nuevo" (D)
$ws.Range("D2").Value = 10000
$ws.Range("D3").Value = 20000
$ws.Range("D4").Value = 10000
$ws.Range("D5").Value = 220202
$ws.Range("D6").Value = 202312
$ws.Range("D7").Value = 254454
$ws.Range("D8").Value = 111111
$ws.Range("D9").Value = 245120
$ws.Range("D10").Value = 2012121
$ws.Range("D11").Value = 12154
$ws.Range("D12").Value = 12124
$ws.Range("D13").Value = 121245
$ws.Range("D14").Value = 121212
$ws.Range("D15").Value = 1212451
$ws.Range("D16").Value = 12541854
$ws.Range("D17").Value = 500001

# Restore the view to the top of the sheet and re-select D18 (matches the
# saved-session selection), clearing the old scrolled-down topLeftCell state.
[void]$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D18").Select()
